$d = $word.ActiveDocument

# 1. "dsn" -> "form.dsn" inside the "Nomor : ... { dsn }" line (the red-colored
#    placeholder run). Using a case-sensitive whole-word search keeps this from
#    touching the unrelated "DSN" run just before it.
$dsnRange = $d.Content
$dsnRange.Find.Execute("dsn", $true, $false, $false, $false, $false, $true, 1, $false, "form.dsn", 2) | Out-Null

# 2. Drop the stray space run right before "penduduk" in the
#    "Nama Lengkap : { penduduk.nama_penduduk}" line, turning
#    "{ penduduk.nama_penduduk}" into "{penduduk.nama_penduduk}".
$spaceRange = $d.Content
$spaceRange.Find.Execute("{ penduduk.nama_penduduk}", $true, $false, $false, $false, $false, $true, 1, $false, "{penduduk.nama_penduduk}", 2) | Out-Null

# 3. Remove the old "_GoBack" bookmark that currently sits between "/" and
#    "penduduks}" in the "Keluarga Yang Datang/Pindah" table.
$oldMark = $d.Bookmarks.Item("_GoBack")
$oldPos = $oldMark.Start
$around = $d.Range($oldPos - 1, $oldPos + 1)
$savedText = $around.Text
$around.Delete()
$d.Range($oldPos - 1, $oldPos - 1).InsertAfter($savedText)

# 4. Re-create "_GoBack" right after the "form.dsn" run we just wrote.
$d.Bookmarks.Add("_GoBack", $d.Range($dsnRange.End, $dsnRange.End))
